$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all Price/Volume cells keep their original text formatting
# (the workbook stores these as plain text, and some values look numeric
# e.g. "1.00", "0.992" -- force text format so Excel does not coerce them).
$ws.Range("D2:E51").NumberFormat = "@"

# --- Cryptocurrency price / volume(1h) refresh ---
$ws.Range("D2").Value = '57.070.08'
$ws.Range("E2").Value = '  +1.12%  '
$ws.Range("D3").Value = '2.399.31'
$ws.Range("E3").Value = '  +2.03%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '504.94'
$ws.Range("E5").Value = '  -1.25%  '
$ws.Range("D6").Value = '132.52'
$ws.Range("E6").Value = '  +4.30%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.37%  '
$ws.Range("D9").Value = '2.413.73'
$ws.Range("E9").Value = '  +2.07%  '
$ws.Range("D10").Value = '0.0967'
$ws.Range("E10").Value = '  +1.36%  '
$ws.Range("E11").Value = '  -1.67%  '
$ws.Range("D12").Value = '0.321'
$ws.Range("E12").Value = '  +1.54%  '
$ws.Range("D13").Value = '4.59'
$ws.Range("E13").Value = '  -4.44%  '
$ws.Range("D14").Value = '2.828.97'
$ws.Range("E14").Value = '  +2.10%  '
$ws.Range("D15").Value = '56.994.55'
$ws.Range("E15").Value = '  +1.14%  '
$ws.Range("D16").Value = '21.80'
$ws.Range("E16").Value = '  +2.02%  '
$ws.Range("E17").Value = '  +2.69%  '
$ws.Range("D18").Value = '2.420.84'
$ws.Range("E18").Value = '  +1.62%  '
$ws.Range("D19").Value = '10.22'
$ws.Range("E19").Value = '  -0.30%  '
$ws.Range("D20").Value = '310.35'
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("D21").Value = '4.03'
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").Value = '6.36'
$ws.Range("E22").Value = '  +4.72%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("E25").Value = '  +0.53%  '
$ws.Range("D26").Value = '0.992'
$ws.Range("E26").Value = '  -0.63%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = '0.153'
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("B28").Value = 'Polygon'
$ws.Range("C28").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D28").Value = '0.377'
$ws.Range("E28").Value = '  -3.06%  '
$ws.Range("D29").Value = '7.47'
$ws.Range("E29").Value = '  +4.34%  '
$ws.Range("D30").Value = '172.73'
$ws.Range("E30").Value = '  -1.05%  '
$ws.Range("E31").Value = '  +1.67%  '
$ws.Range("E32").Value = '  +0.41%  '
$ws.Range("D33").Value = '1.13'
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("D34").Value = '5.93'
$ws.Range("E34").Value = '  -2.86%  '
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("D37").Value = '17.95'
$ws.Range("E37").Value = '  +2.11%  '
$ws.Range("D38").Value = '1.20'
$ws.Range("E38").Value = '  +1.32%  '
$ws.Range("D39").Value = '3.83'
$ws.Range("E39").Value = '  +3.55%  '
$ws.Range("D40").Value = '36.65'
$ws.Range("E40").Value = '  +3.49%  '
$ws.Range("D41").Value = '0.805'
$ws.Range("E41").Value = '  +0.53%  '
$ws.Range("E42").Value = '  +1.36%  '
$ws.Range("D43").Value = '132.40'
$ws.Range("E43").Value = '  +9.75%  '
$ws.Range("D44").Value = '4.86'
$ws.Range("E44").Value = '  +0.85%  '
$ws.Range("D45").Value = '3.35'
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").Value = '0.569'
$ws.Range("E46").Value = '  +0.63%  '
$ws.Range("D47").Value = '252.05'
$ws.Range("E47").Value = '  -0.27%  '
$ws.Range("D48").Value = '0.0910'
$ws.Range("E48").Value = '  +0.71%  '
$ws.Range("D49").Value = '0.0487'
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = '0.0209'
$ws.Range("E50").Value = '  +1.15%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '16.96'
$ws.Range("E51").Value = '  +2.45%  '
